$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$origStyle = $cell.Style

# B11 currently contains the text "R40"; change it to the text "1"
$cell.Value = "'1"
$cell.Style = $origStyle

$wb.Save()
